$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 28 de Abril de 2020 a las 11:22"

# Row 35 - Rumania: update Recuperados, Muertes hoy, Muertes
$ws.Range("E35").Value = 7548
$ws.Range("G35").Value = 9
$ws.Range("H35").Value = 650

# Rows 39/40 swap places (Ucrania <-> Indonesia) with updated figures
# Row 39 becomes Indonesia with new data
$ws.Range("A39").Value = "Indonesia"
$ws.Range("B39").Value = 9511
$ws.Range("C39").Value = 415
$ws.Range("D39").Value = 1254
$ws.Range("E39").Value = 7484
$ws.Range("F39").Value = 0
$ws.Range("G39").Value = 8
$ws.Range("H39").Value = 773

# Row 40 becomes Ucrania, carrying the former row 39 figures
$ws.Range("A40").Value = "Ucrania"
$ws.Range("B40").Value = 9410
$ws.Range("C40").Value = 401
$ws.Range("D40").Value = 992
$ws.Range("E40").Value = 8179
$ws.Range("F40").Value = 121
$ws.Range("G40").Value = 19
$ws.Range("H40").Value = 239

# Row 46 - Australia: updated figures
$ws.Range("B46").Value = 6731
$ws.Range("C46").Value = 11
$ws.Range("D46").Value = 5626
$ws.Range("E46").Value = 1021
$ws.Range("F46").Value = 42

# Rows 47/48/49 shift (Banglades moves up ahead of Republica Dominicana/Panama)
# Row 47 becomes Banglades with new data
$ws.Range("A47").Value = "Banglades"
$ws.Range("B47").Value = 6462
$ws.Range("C47").Value = 549
$ws.Range("D47").Value = 139
$ws.Range("E47").Value = 6168
$ws.Range("F47").Value = 1
$ws.Range("G47").Value = 3
$ws.Range("H47").Value = 155

# Row 48 becomes Republica Dominicana, carrying the former row 47 figures
$ws.Range("A48").Value = "Republica Dominicana"
$ws.Range("B48").Value = 6293
$ws.Range("C48").Value = 0
$ws.Range("D48").Value = 993
$ws.Range("E48").Value = 5018
$ws.Range("F48").Value = 144
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 282

# Row 49 becomes Panama, carrying the former row 48 figures
$ws.Range("A49").Value = "Panama"
$ws.Range("B49").Value = 6021
$ws.Range("C49").Value = 242
$ws.Range("D49").Value = 455
$ws.Range("E49").Value = 5399
$ws.Range("F49").Value = 89
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = 167

# Row 68 - Uzbekistan: updated figures
$ws.Range("B68").Value = 1939
$ws.Range("C68").Value = 35
$ws.Range("D68").Value = 909
$ws.Range("E68").Value = 1022
